# Generate Report for Handback
# Update status/timestamp values for the a8519961-... record across the
# Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Row 3 (a8519961...) and Row 4 (c5d600b8...) share the same
# "Latest HO Xliff Generate Date" value; both move forward together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 16:14:52"
$wsOverview.Range("G4").Value = "2016-08-24 16:14:52"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority: ht -> mt for both the a8519961 and c5d600b8 rows
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-08-24 16:14:47"
$wsZhCn.Range("H4").Value = "2016-08-24 16:14:47"
# Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-08-24 16:15:25"
$wsZhCn.Range("K4").Value = "2016-08-24 16:15:25"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority: ht -> mt for both the a8519961 and c5d600b8 rows
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime (shares value with Overview's G3/G4)
$wsDeDe.Range("H3").Value = "2016-08-24 16:14:52"
$wsDeDe.Range("H4").Value = "2016-08-24 16:14:52"
# Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-24 16:15:33"
$wsDeDe.Range("K4").Value = "2016-08-24 16:15:33"
